$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.57369065284729
$ws.Range("B1").Value = 1.359740257263184
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.699376583099365
$ws.Range("E1").Value = 1.479905128479004
